# Weekly update: a new price record for "Haba" (Femacal de La Calera) is
# inserted as row 296, pushing the existing rows 296-321 down to 297-322.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 296 (shifts 296:321 -> 297:322,
# carrying over the date-format style already used by column D).
$ws.Rows(296).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A296").Value = 3
$ws.Range("B296").Value = "Femacal de La Calera"
$ws.Range("C296").Value = "Coquimbo"
$ws.Range("D296").Value = "2023-10-24"
$ws.Range("E296").Value = 5
$ws.Range("F296").Value = 100112026
$ws.Range("G296").Value = "Haba"
$ws.Range("H296").Value = "Sin especificar"
$ws.Range("I296").Value = "Primera"
$ws.Range("J296").Value = 55
$ws.Range("K296").Value = 9000
$ws.Range("L296").Value = 9000
$ws.Range("M296").Value = 9000
$ws.Range("N296").Value = '$/saco 25 kilos'
$ws.Range("O296").Value = "Provincia de Petorca"
$ws.Range("P296").Value = 360
$ws.Range("Q296").Value = 25
$ws.Range("R296").Value = "Hortaliza"
